$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.058286333333333
$ws.Cells.Item(2, 8).Value = 3.174859
$ws.Cells.Item(2, 9).Value = 0.6553662741588026
$ws.Cells.Item(2, 10).Value = 0.6553662741588027
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 101.7913436666667
$ws.Cells.Item(2, 14).Value = 305.374031
$ws.Cells.Item(2, 15).Value = 0.2805454758424659
$ws.Cells.Item(2, 16).Value = 0.2805454758424659
$ws.Cells.Item(2, 17).Value = 107.7243878540699
$ws.Cells.Item(2, 18).Value = 969.519490686629
$ws.Cells.Item(2, 19).Value = 0.1838600432349852
$ws.Cells.Item(2, 20).Value = 0.1838600432349853

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.058286333333333
$ws.Cells.Item(3, 8).Value = 3.174859
$ws.Cells.Item(3, 9).Value = 0.6553662741588026
$ws.Cells.Item(3, 10).Value = 0.6553662741588027
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 54.34621166666667
$ws.Cells.Item(3, 14).Value = 163.038635
$ws.Cells.Item(3, 15).Value = 0.1497827149446808
$ws.Cells.Item(3, 16).Value = 0.1497827149446808
$ws.Cells.Item(3, 17).Value = 57.51385307527389
$ws.Cells.Item(3, 18).Value = 517.6246776774649
$ws.Cells.Item(3, 19).Value = 0.09816253982668546
$ws.Cells.Item(3, 20).Value = 0.09816253982668548

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.058286333333333
$ws.Cells.Item(4, 8).Value = 3.174859
$ws.Cells.Item(4, 9).Value = 0.6553662741588026
$ws.Cells.Item(4, 10).Value = 0.6553662741588027
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 122.2542826666667
$ws.Cells.Item(4, 14).Value = 366.762848
$ws.Cells.Item(4, 15).Value = 0.3369430510399163
$ws.Cells.Item(4, 16).Value = 0.3369430510399163
$ws.Cells.Item(4, 17).Value = 129.3800365376035
$ws.Cells.Item(4, 18).Value = 1164.420328838432
$ws.Cells.Item(4, 19).Value = 0.2208211119637292
$ws.Cells.Item(4, 20).Value = 0.2208211119637292

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.058286333333333
$ws.Cells.Item(5, 8).Value = 3.174859
$ws.Cells.Item(5, 9).Value = 0.6553662741588026
$ws.Cells.Item(5, 10).Value = 0.6553662741588027
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 75.13712199999999
$ws.Cells.Item(5, 14).Value = 225.411366
$ws.Cells.Item(5, 15).Value = 0.2070842066291165
$ws.Cells.Item(5, 16).Value = 0.2070842066291166
$ws.Cells.Item(5, 17).Value = 79.51658933859932
$ws.Cells.Item(5, 18).Value = 715.6493040473939
$ws.Cells.Item(5, 19).Value = 0.1357160049356557
$ws.Cells.Item(5, 20).Value = 0.1357160049356557

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.058286333333333
$ws.Cells.Item(6, 8).Value = 3.174859
$ws.Cells.Item(6, 9).Value = 0.6553662741588026
$ws.Cells.Item(6, 10).Value = 0.6553662741588027
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 9.304706666666666
$ws.Cells.Item(6, 14).Value = 27.91412
$ws.Cells.Item(6, 15).Value = 0.02564455154382035
$ws.Cells.Item(6, 16).Value = 0.02564455154382035
$ws.Cells.Item(6, 17).Value = 9.847043901008888
$ws.Cells.Item(6, 18).Value = 88.62339510907998
$ws.Cells.Item(6, 19).Value = 0.01680657419774691
$ws.Cells.Item(6, 20).Value = 0.01680657419774691

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.556515
$ws.Cells.Item(7, 8).Value = 1.669545
$ws.Cells.Item(7, 9).Value = 0.3446337258411974
$ws.Cells.Item(7, 10).Value = 0.3446337258411974
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 101.7913436666667
$ws.Cells.Item(7, 14).Value = 305.374031
$ws.Cells.Item(7, 15).Value = 0.2805454758424659
$ws.Cells.Item(7, 16).Value = 0.2805454758424659
$ws.Cells.Item(7, 17).Value = 56.64840962065499
$ws.Cells.Item(7, 18).Value = 509.835686585895
$ws.Cells.Item(7, 19).Value = 0.09668543260748066
$ws.Cells.Item(7, 20).Value = 0.09668543260748068

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.556515
$ws.Cells.Item(8, 8).Value = 1.669545
$ws.Cells.Item(8, 9).Value = 0.3446337258411974
$ws.Cells.Item(8, 10).Value = 0.3446337258411974
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 54.34621166666667
$ws.Cells.Item(8, 14).Value = 163.038635
$ws.Cells.Item(8, 15).Value = 0.1497827149446808
$ws.Cells.Item(8, 16).Value = 0.1497827149446808
$ws.Cells.Item(8, 17).Value = 30.244481985675
$ws.Cells.Item(8, 18).Value = 272.200337871075
$ws.Cells.Item(8, 19).Value = 0.05162017511799535
$ws.Cells.Item(8, 20).Value = 0.05162017511799535

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.556515
$ws.Cells.Item(9, 8).Value = 1.669545
$ws.Cells.Item(9, 9).Value = 0.3446337258411974
$ws.Cells.Item(9, 10).Value = 0.3446337258411974
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 122.2542826666667
$ws.Cells.Item(9, 14).Value = 366.762848
$ws.Cells.Item(9, 15).Value = 0.3369430510399163
$ws.Cells.Item(9, 16).Value = 0.3369430510399163
$ws.Cells.Item(9, 17).Value = 68.03634211823999
$ws.Cells.Item(9, 18).Value = 612.3270790641599
$ws.Cells.Item(9, 19).Value = 0.1161219390761871
$ws.Cells.Item(9, 20).Value = 0.1161219390761871

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.556515
$ws.Cells.Item(10, 8).Value = 1.669545
$ws.Cells.Item(10, 9).Value = 0.3446337258411974
$ws.Cells.Item(10, 10).Value = 0.3446337258411974
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 75.13712199999999
$ws.Cells.Item(10, 14).Value = 225.411366
$ws.Cells.Item(10, 15).Value = 0.2070842066291165
$ws.Cells.Item(10, 16).Value = 0.2070842066291166
$ws.Cells.Item(10, 17).Value = 41.81493544983
$ws.Cells.Item(10, 18).Value = 376.33441904847
$ws.Cells.Item(10, 19).Value = 0.07136820169346082
$ws.Cells.Item(10, 20).Value = 0.07136820169346084

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.556515
$ws.Cells.Item(11, 8).Value = 1.669545
$ws.Cells.Item(11, 9).Value = 0.3446337258411974
$ws.Cells.Item(11, 10).Value = 0.3446337258411974
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 9.304706666666666
$ws.Cells.Item(11, 14).Value = 27.91412
$ws.Cells.Item(11, 15).Value = 0.02564455154382035
$ws.Cells.Item(11, 16).Value = 0.02564455154382035
$ws.Cells.Item(11, 17).Value = 5.1782088306
$ws.Cells.Item(11, 18).Value = 46.60387947539999
$ws.Cells.Item(11, 19).Value = 0.00883797734607344
$ws.Cells.Item(11, 20).Value = 0.00883797734607344
